$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.076.16"
$ws.Range("E2").Value = "  +6.13%  "
$ws.Range("D3").Value = "3.553.99"
$ws.Range("E3").Value = "  +9.79%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "552.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("D7").Value = "3.548.42"
$ws.Range("E7").Value = "  +9.53%  "
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.633"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").Value = "4.126.91"
$ws.Range("E15").Value = "  +10.09%  "
$ws.Range("D16").Value = "3.557.67"
$ws.Range("E16").Value = "  +10.25%  "
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "67.130.38"
$ws.Range("E18").Value = "  +6.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "433.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  +7.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "646.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  +5.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +23.16%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0824"
$ws.Range("E38").Value = "  +15.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.08%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("D45").Value = "3.034.48"
$ws.Range("E45").Value = "  +5.22%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +13.09%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0419"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.73%  "
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.131"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.00%  "

Write-Host "Applied all changes"